# Apply the cryptos-list price/volume refresh described by the commit diff.
# Values are read from the live coinranking.com feed on each scheduled run;
# this script pins the specific before -> after values captured in that run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new text would otherwise be auto-parsed as a number by
# Excel (e.g. "1.001", "0.7500") must be forced to Text so the stored value
# keeps its exact literal form (and the "1,001" style thousands grouping in
# the source data is preserved as text, not coerced to a float).
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D20", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.933.10"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "1.876.84"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "0.7390"
$ws.Range("E5").Value = "  -3.69%  "

$ws.Range("D6").Value = "242.95"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "0.3149"
$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("D9").Value = "0.07239"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").Value = "24.67"
$ws.Range("E10").Value = "  -3.92%  "

$ws.Range("D11").Value = "0.08346"
$ws.Range("E11").Value = "  -2.26%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.919.60"
$ws.Range("E12").Value = "  +1.85%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7500"
$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("D14").Value = "5.391"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").Value = "92.32"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").Value = "29.935.20"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").Value = "6.107"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "246.78"
$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("E19").Value = "  -1.46%  "

$ws.Range("D20").Value = "0.000007843"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "2.141.75"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").Value = "8.006"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D26").Value = "9.305"
$ws.Range("E26").Value = "  -1.27%  "

$ws.Range("D27").Value = "165.11"
$ws.Range("E27").Value = "  +1.37%  "

$ws.Range("D28").Value = "18.66"
$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("D29").Value = "2.022"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").Value = "1.497"
$ws.Range("E30").Value = "  +2.50%  "

$ws.Range("D31").Value = "4.613"
$ws.Range("E31").Value = "  +2.45%  "

$ws.Range("D32").Value = "1.538"
$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("D33").Value = "4.254"
$ws.Range("E33").Value = "  +3.86%  "

$ws.Range("D34").Value = "0.05340"
$ws.Range("E34").Value = "  -2.02%  "

$ws.Range("D35").Value = "1.235"
$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("D36").Value = "0.7500"
$ws.Range("E36").Value = "  +0.91%  "

$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").Value = "2.702"
$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").Value = "0.01960"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").Value = "2.755"
$ws.Range("E40").Value = "  -0.95%  "

$ws.Range("D41").Value = "0.4523"
$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("D42").Value = "1.114.12"
$ws.Range("E42").Value = "  +1.10%  "

$ws.Range("D43").Value = "6.135"

$ws.Range("D44").Value = "72.49"
$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("D45").Value = "0.8627"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("D46").Value = "104.59"
$ws.Range("E46").Value = "  +1.56%  "

$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").Value = "1.864"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "7.598"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("D50").Value = "9.510"
$ws.Range("E50").Value = "  -2.30%  "

$ws.Range("D51").Value = "2.037.54"
$ws.Range("E51").Value = "  -0.07%  "
